# Change the color of school points to yellow
# (Underlying data refresh: regenerate the child/school route dataset.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$NBSP = [char]0x00A0

# Helper: write a value while making sure purely-numeric-looking strings
# (e.g. "8", "46.0") are stored as TEXT, not converted to a Number cell,
# and without leaving a lingering style on the cell.
function Set-TextValue {
    param($range, $value)

    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- nChildren ---
Set-TextValue $ws.Range("B4") "8"

# --- child row 6 (index 0) ---
Set-TextValue $ws.Range("B6") "2"
$ws.Range("C6").Value = "Elwanda " + $NBSP
$ws.Range("D6").Value = "Cassy " + $NBSP
$ws.Range("E6").Value = "-5.57,-9.66"
$ws.Range("F6").Value = "Tamisha(mother): 0550693864"
Set-TextValue $ws.Range("H6") "46.0"

# --- child row 7 (index 1) ---
Set-TextValue $ws.Range("B7") "20"
$ws.Range("C7").Value = "Ron"
$ws.Range("D7").Value = "Cohen"
$ws.Range("E7").Value = "-8.77,-7.22"
$ws.Range("F7").Value = "Bernardine(mother): 0576270618"
Set-TextValue $ws.Range("H7") "41.0"

# --- child row 8 (index 2) ---
Set-TextValue $ws.Range("B8") "6"
$ws.Range("C8").Value = "Ema " + $NBSP
$ws.Range("D8").Value = "Ardell " + $NBSP
$ws.Range("E8").Value = "-3.09,-3.85"
$ws.Range("F8").Value = "Carley(grandmother): 0533587167"
$ws.Range("G8").Value = "7:13:00"
Set-TextValue $ws.Range("H8") "33.0"

# --- child row 9 (index 3) ---
Set-TextValue $ws.Range("B9") "3"
$ws.Range("C9").Value = "Alexia " + $NBSP
$ws.Range("D9").Value = "Ramonita " + $NBSP
$ws.Range("E9").Value = "-1.78,-3.65"
$ws.Range("F9").Value = "Han(father): 0567537032"
$ws.Range("G9").Value = "7:15:00"
Set-TextValue $ws.Range("H9") "31.0"

# --- child row 10 (index 4) ---
Set-TextValue $ws.Range("B10") "5"
$ws.Range("C10").Value = "Patti " + $NBSP
$ws.Range("D10").Value = "Lavenia " + $NBSP
$ws.Range("E10").Value = "-0.56,-6.13"
$ws.Range("F10").Value = "Jennell(mother): 0503029941"
$ws.Range("G10").Value = "7:19:00"
Set-TextValue $ws.Range("H10") "27.0"

# --- child row 11 (index 5) ---
Set-TextValue $ws.Range("B11") "18"
$ws.Range("C11").Value = "Kandis " + $NBSP
$ws.Range("D11").Value = "Zulma " + $NBSP
$ws.Range("E11").Value = "1.22,-9.32"
$ws.Range("F11").Value = "Kylie(mother): 0575413269"
$ws.Range("G11").Value = "7:24:00"
Set-TextValue $ws.Range("H11") "22.0"

# --- child row 12 (index 6) ---
Set-TextValue $ws.Range("B12") "4"
$ws.Range("C12").Value = "Francisca " + $NBSP
$ws.Range("D12").Value = "Stevie " + $NBSP
$ws.Range("E12").Value = "7.54,-8.77"
$ws.Range("F12").Value = "Bernardine(mother): 0561339273"
$ws.Range("G12").Value = "7:31:00"
Set-TextValue $ws.Range("H12") "15.0"

# --- child row 13 (index 7) -- newly added child ---
Set-TextValue $ws.Range("A13") "7"
Set-TextValue $ws.Range("B13") "10"
$ws.Range("C13").Value = "Demetra " + $NBSP
$ws.Range("D13").Value = "Francene " + $NBSP
$ws.Range("E13").Value = "1.17,-4.35"
$ws.Range("F13").Value = "Dorian(mother): 0534328089"
$ws.Range("G13").Value = "7:39:00"
Set-TextValue $ws.Range("H13") "7.0"

# --- school row, now shifted to row 14 ---
$ws.Range("A14").Value = "school"
Set-TextValue $ws.Range("B14") "3"
$ws.Range("C14").Value = "Ironiah"
$ws.Range("D14").Value = "mySchool"
$ws.Range("E14").Value = "0,0"
$ws.Range("F14").Value = "Shir(secretary): 0523345098"
$ws.Range("G14").Value = "7:46:00"

# --- cost row, now shifted to row 15 ---
$ws.Range("A15").Value = "cost"
Set-TextValue $ws.Range("B15") "25"

# --- time row, now shifted to row 16 ---
$ws.Range("A16").Value = "time"
Set-TextValue $ws.Range("B16") "46.0"
